$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 / C10: Angrist & Piscke (MHE) reading now leads with a bullet and gains
#     the "Topic and 3 Questions Due" note that used to live on row 8 ---
# (Written first so the new shared-string entries land in the same order as the target
#  workbook: "Angrist ... MHE" before "Imbens and Rubin".)
$ws.Range("C10").Formula = "'- Angrist and Piscke (MHE) Chapter 1 (pages 3-24)`n- **Topic and 3 Questions Due**"

# Bold run: "ages 3-24)" (chars 40-49)
$c10run2 = $ws.Range("C10").Characters(40, 10)
$c10run2.Font.Name = "SFBX1200"

# Regular run: "\n- **Topic and 3 Questions Due**" (chars 50-81)
$c10run3 = $ws.Range("C10").Characters(50, 32)
$c10run3.Font.Name = "SFRM1200"

$ws.Rows(10).RowHeight = 51

# --- Row 8 / C8: "Imbens and Rubin" reading replaces "Morgan and Winship" reading ---
# Leading apostrophe forces Excel's literal-text (quote-prefix) handling, matching the
# "s=11" quotePrefix style seen in the target workbook for this cell.
$ws.Range("C8").Formula = "'- Imbens and Rubin. Chapter 1 (pages 3-22).`n(PDF on Sakai for Duke students)"

# Bold run: "ages 3-22)" (chars 33-42)
$c8run2 = $ws.Range("C8").Characters(33, 10)
$c8run2.Font.Name = "SFBX1200"

# Regular run: ".\n(PDF on Sakai for Duke students)" (chars 43-76)
$c8run3 = $ws.Range("C8").Characters(43, 34)
$c8run3.Font.Name = "SFRM1200"

$ws.Rows(8).RowHeight = 68

# --- Sheet view / selection moved down as the author scrolled to the edited rows ---
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("C9").Select()
